# Update MSME Country Indicators - Ecuador Summary figures to more precise
# values. The target cells already hold their numbers as text (shared
# strings rather than numeric cells), so assign the new values with a
# leading apostrophe - same as typing '40.76 into Excel - which keeps the
# cell text-typed instead of letting Excel auto-convert the numeric-looking
# string into a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B11" = "40.76"
    "C11" = "4.49"
    "D11" = "45.25"
    "B12" = "21.38"
    "C12" = "39.15"
    "D12" = "60.53"
    "B14" = "89.62"
    "C14" = "9.87"
    "B30" = "31.65"
    "C30" = "1.46"
    "D30" = "33.11"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
